$d = $word.ActiveDocument

# Update the header date paragraph
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("2024-01-03 Wednesday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-01-04 Thursday", 2) | Out-Null

# Update each math-expression cell in the 20x5 table, in row-major order,
# matching the document/diff order exactly (no text search, so no risk of
# accidental substring collisions between old/new expressions).
$newValues = @(
    '49+20=',
    '70+12=',
    '38-7=',
    '42-0=',
    '39-21=',
    '7+47=',
    '59+29=',
    '98-4=',
    '32+47=',
    '27+66=',
    '18-9=',
    '80+17=',
    '9+75=',
    '8+22=',
    '62+6=',
    '31+42=',
    '64+9=',
    '48+47=',
    '41-26=',
    '26+62=',
    '97-25=',
    '5+9=',
    '50+1=',
    '64+12=',
    '93-90=',
    '69+28=',
    '80+11=',
    '80+18=',
    '48+24=',
    '96-78=',
    '86-57=',
    '76+5=',
    '17+17=',
    '56-54=',
    '54+8=',
    '51-16=',
    '24+53=',
    '40+58=',
    '86-40=',
    '11+28=',
    '44-37=',
    '18+49=',
    '71-42=',
    '90-77=',
    '38+11=',
    '26+33=',
    '10+1=',
    '43+6=',
    '63+3=',
    '55+6=',
    '57+4=',
    '28+34=',
    '99-51=',
    '9+64=',
    '71-60=',
    '34+13=',
    '62-61=',
    '36+28=',
    '83-40=',
    '75-12=',
    '67-54=',
    '19+67=',
    '62-56=',
    '39-20=',
    '1+40=',
    '81-80=',
    '52+21=',
    '62-30=',
    '42+31=',
    '66+31=',
    '19+76=',
    '56+35=',
    '77+0=',
    '96-60=',
    '80+9=',
    '26+36=',
    '77-7=',
    '12-11=',
    '38-28=',
    '1+96=',
    '62+12=',
    '29-3=',
    '98-84=',
    '0+31=',
    '85-66=',
    '67-22=',
    '83-62=',
    '93-70=',
    '6+63=',
    '88-11=',
    '50+37=',
    '68+7=',
    '53+36=',
    '98-1=',
    '86-84=',
    '68-53=',
    '76-71=',
    '77-69=',
    '66-62=',
    '22+10=',
)

$t = $d.Tables(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cellRange = $cell.Range
        $cellRange.SetRange($cellRange.Start, $cellRange.End - 1)
        $cellRange.Text = $newValues[$i]
        $i = $i + 1
    }
}

Write-Output "done: $i cells updated"
